$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 1.74
$ws.Range("L2").Value = 1.25
$ws.Range("N2").Value = 6.4
$ws.Range("Q2").Value = 1.5
$ws.Range("S2").Value = 2.2
$ws.Range("T2").Value = 1.54
$ws.Range("U2").Value = 2.66
$ws.Range("AB2").Value = 15
$ws.Range("AF2").Value = 15

# Row 3
$ws.Range("F3").Value = 1.93
$ws.Range("H3").Value = 4.1
$ws.Range("J3").Value = 3.95
$ws.Range("K3").Value = 4.1
$ws.Range("P3").Value = 2.28
$ws.Range("R3").Value = 1.52
$ws.Range("S3").Value = 2.68
$ws.Range("U3").Value = 2.32
$ws.Range("X3").Value = 23
$ws.Range("Y3").Value = 19.5
$ws.Range("AA3").Value = 95
$ws.Range("AI3").Value = 48
$ws.Range("AL3").Value = 30
$ws.Range("AM3").Value = 75
$ws.Range("AO3").Value = 38

# Row 4
$ws.Range("F4").Value = 1.02
$ws.Range("H4").Value = 1.02
$ws.Range("J4").Value = 1.02
$ws.Range("L4").Value = 1.01
$ws.Range("M4").Value = 1.01
$ws.Range("N4").Value = 1.17
$ws.Range("O4").Value = 1.01
$ws.Range("P4").Value = 1.17
$ws.Range("R4").Value = 1.09
$ws.Range("S4").Value = 1.39
$ws.Range("T4").Value = 1.04
$ws.Range("U4").Value = 1.04
$ws.Range("V4").Value = 1.01
$ws.Range("W4").Value = 1.01
$ws.Range("X4").Value = 990
$ws.Range("Y4").Value = 990
$ws.Range("Z4").Value = 1000
$ws.Range("AA4").Value = 1000
$ws.Range("AB4").Value = 990
$ws.Range("AC4").Value = 990
$ws.Range("AD4").Value = 990
$ws.Range("AE4").Value = 1000
$ws.Range("AF4").Value = 1000
$ws.Range("AG4").Value = 990
$ws.Range("AH4").Value = 990
$ws.Range("AI4").Value = 1000
$ws.Range("AJ4").Value = 1000
$ws.Range("AK4").Value = 1000
$ws.Range("AL4").Value = 1000
$ws.Range("AM4").Value = 1000
$ws.Range("AN4").Value = 1000
$ws.Range("AO4").Value = 1000

# Row 6
$ws.Range("F6").Value = 2.28
$ws.Range("J6").Value = 3.45
$ws.Range("K6").Value = 3.5
$ws.Range("N6").Value = 3.55
$ws.Range("O6").Value = 1.37
$ws.Range("P6").Value = 1.83
$ws.Range("Q6").Value = 2.16
$ws.Range("R6").Value = 1.33
$ws.Range("S6").Value = 3.85
$ws.Range("T6").Value = 1.87
$ws.Range("U6").Value = 2.08
$ws.Range("X6").Value = 13
$ws.Range("Y6").Value = 13
$ws.Range("AB6").Value = 9.4
$ws.Range("AF6").Value = 14
$ws.Range("AH6").Value = 19.5
$ws.Range("AJ6").Value = 32
$ws.Range("AK6").Value = 26
$ws.Range("AM6").Value = 130
$ws.Range("AN6").Value = 21

# Row 7
$ws.Range("F7").Value = 1.58
$ws.Range("G7").Value = 1.59
$ws.Range("H7").Value = 6.4
$ws.Range("I7").Value = 6.6
$ws.Range("N7").Value = 5.2
$ws.Range("O7").Value = 1.23
$ws.Range("P7").Value = 2.44
$ws.Range("Q7").Value = 1.68
$ws.Range("R7").Value = 1.55
$ws.Range("S7").Value = 2.74
$ws.Range("U7").Value = 2.24
$ws.Range("Y7").Value = 27
$ws.Range("AA7").Value = 200
$ws.Range("AG7").Value = 9.6
$ws.Range("AH7").Value = 21
$ws.Range("AJ7").Value = 15
$ws.Range("AK7").Value = 14.5
$ws.Range("AN7").Value = 7
$ws.Range("AO7").Value = 85

# Row 8
$ws.Range("F8").Value = 1.91
$ws.Range("G8").Value = 2.12
$ws.Range("H8").Value = 3.45
$ws.Range("I8").Value = 4.5
$ws.Range("J8").Value = 3.95
$ws.Range("K8").Value = 4.7
$ws.Range("P8").Value = 2.5
$ws.Range("Q8").Value = 1.52

# Row 9
$ws.Range("F9").Value = 2.3
$ws.Range("G9").Value = 2.32
$ws.Range("H9").Value = 3.35
$ws.Range("I9").Value = 3.45
$ws.Range("N9").Value = 3.95
$ws.Range("O9").Value = 1.32
$ws.Range("P9").Value = 1.97
$ws.Range("Q9").Value = 2
$ws.Range("R9").Value = 1.38
$ws.Range("S9").Value = 3.45
$ws.Range("T9").Value = 1.79
$ws.Range("U9").Value = 2.2
$ws.Range("AA9").Value = 70
$ws.Range("AB9").Value = 10.5
$ws.Range("AC9").Value = 7.8
$ws.Range("AD9").Value = 14.5
$ws.Range("AE9").Value = 38
$ws.Range("AF9").Value = 15
$ws.Range("AI9").Value = 60
$ws.Range("AL9").Value = 44
$ws.Range("AN9").Value = 18
$ws.Range("AO9").Value = 36

# Row 10
$ws.Range("H10").Value = 1.46
$ws.Range("I10").Value = 1.47
$ws.Range("N10").Value = 5
$ws.Range("O10").Value = 1.23
$ws.Range("P10").Value = 2.36
$ws.Range("R10").Value = 1.54
$ws.Range("S10").Value = 2.76
$ws.Range("Y10").Value = 9.800000000000001
$ws.Range("Z10").Value = 8.800000000000001
$ws.Range("AA10").Value = 13
$ws.Range("AB10").Value = 28
$ws.Range("AH10").Value = 23
$ws.Range("AO10").Value = 6.4
